# Implemented All planned Methods in Data Access Layer
#
# This mirrors the author's commit: every method row's "Implementation"
# status (column D) is moved from "Not Started" to "Done" (the ones that
# were already "Done" stay "Done"), the now-finished method
# "DoesPersonHaveUser44" is renamed to "DoesPersonHaveUser" (the trailing
# "44" was a placeholder/typo), and the view is scrolled back towards the
# top-left of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")

# --- Column D ("Implementation") -> "Done" for every method row ---------
$ws.Range("D3:D7").Value = "Done"
$ws.Range("D9:D18").Value = "Done"
$ws.Range("D20:D22").Value = "Done"
$ws.Range("D24:D30").Value = "Done"
$ws.Range("D32:D36").Value = "Done"
$ws.Range("D38:D43").Value = "Done"
$ws.Range("D45:D51").Value = "Done"
$ws.Range("D53:D62").Value = "Done"
$ws.Range("D64:D71").Value = "Done"
$ws.Range("D73:D78").Value = "Done"
$ws.Range("D80:D86").Value = "Done"
$ws.Range("D88:D91").Value = "Done"
$ws.Range("D93:D104").Value = "Done"
$ws.Range("D106:D110").Value = "Done"

# --- Rename the finished method (drop the placeholder "44") -------------
$ws.Range("C103").Value = "public static bool DoesPersonHaveUser(int PersonID)"

# --- Recalculate so the progress-percentage formulas (I4/J4) refresh ----
$excel.Calculate()

# --- Restore the view to the top-left of the sheet -----------------------
$ws.Range("F9").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
